$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D6 currently stores the text "4:36" (a shared string). The author
# replaced it with the plain number 436. Writing a numeric value here
# converts the cell from a shared-string cell to a numeric cell, and the
# now-unused "4:36" shared-string entry is dropped on save (shifting the
# later "17:05" entry's index down by one, e.g. E6 which references it).
$ws.Range("D6").Value = 436
